$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Alias" value for Augusta Delono (row 4) changes from "AUD" to "ADO",
# and the old value "AUD" is copied into the new cell J4.
$ws.Range("J4").Value = $ws.Range("D4").Value()
$ws.Range("D4").Value = "ADO"

$ws.Range("H17").Select()
